$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B2 text change and C2 value change
$ws.Range("B2").Value = "<may>"
$ws.Range("C2").Value = 15

# Row 3: C3 value change
$ws.Range("C3").Value = 15

# Row 5: C5 value change
$ws.Range("C5").Value = 17

# Row 6: C6 value change
$ws.Range("C6").Value = 14

# Row 7: C7 value change
$ws.Range("C7").Value = 10

# Row 8: C8 value change
$ws.Range("C8").Value = 10

# Row 9: B9 text change
$ws.Range("B9").Value = "<by>"

# Row 10: B10 text change and C10 value change
$ws.Range("B10").Value = "<alt>"
$ws.Range("C10").Value = 10

# Row 11: C11 value change
$ws.Range("C11").Value = 15

# Row 13: C13 value change
$ws.Range("C13").Value = 13

# Row 14: C14 value change
$ws.Range("C14").Value = 11

# Row 15: C15 value change
$ws.Range("C15").Value = 17

# Row 16: C16 value change
$ws.Range("C16").Value = 16

# Row 17: C17 value change
$ws.Range("C17").Value = 17

# Row 18: C18 value change
$ws.Range("C18").Value = 10
